# Weekly update: insert 3 new report rows (week of 2022-06-02) above the
# two most recent existing weeks, pushing those rows down.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new rows before current row 71 (shifts old rows 71-74 -> 74-77,
# carrying their existing formatting/content down with them).
$ws.Rows("71:73").Insert()

# Copy the date cell formatting (style) used throughout column D into the
# newly inserted rows so the date values render/format correctly.
$ws.Range("D74").Copy()
$ws.Range("D71:D73").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 71: Especial quality, 2022-06-02, Región de O'Higgins
$ws.Range("A71").Value = 2
$ws.Range("B71").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C71").Value = "Coquimbo"
$ws.Range("D71").Value = 44714
$ws.Range("E71").Value = 4
$ws.Range("F71").Value = "Fruta"
$ws.Range("G71").Value = 100103
$ws.Range("H71").Value = "Frutos de hueso (carozo)"
$ws.Range("I71").Value = 100103002
$ws.Range("J71").Value = "Ciruela"
$ws.Range("K71").Value = "Angeleno"
$ws.Range("L71").Value = "Especial"
$ws.Range("M71").Value = 10
$ws.Range("N71").Value = 220000
$ws.Range("O71").Value = 230000
$ws.Range("P71").Value = 225000
$ws.Range("Q71").Value = "$/bins (450 kilos)"
$ws.Range("R71").Value = "Región de O'Higgins"
$ws.Range("S71").Value = 500
$ws.Range("T71").Value = 450

# Row 72: Primera quality, 2022-06-02, Región de O'Higgins
$ws.Range("A72").Value = 2
$ws.Range("B72").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C72").Value = "Coquimbo"
$ws.Range("D72").Value = 44714
$ws.Range("E72").Value = 4
$ws.Range("F72").Value = "Fruta"
$ws.Range("G72").Value = 100103
$ws.Range("H72").Value = "Frutos de hueso (carozo)"
$ws.Range("I72").Value = 100103002
$ws.Range("J72").Value = "Ciruela"
$ws.Range("K72").Value = "Angeleno"
$ws.Range("L72").Value = "Primera"
$ws.Range("M72").Value = 16
$ws.Range("N72").Value = 200000
$ws.Range("O72").Value = 210000
$ws.Range("P72").Value = 205000
$ws.Range("Q72").Value = "$/bins (450 kilos)"
$ws.Range("R72").Value = "Región de O'Higgins"
$ws.Range("S72").Value = 456
$ws.Range("T72").Value = 450

# Row 73: Segunda quality, 2022-06-02, Región de O'Higgins
$ws.Range("A73").Value = 2
$ws.Range("B73").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C73").Value = "Coquimbo"
$ws.Range("D73").Value = 44714
$ws.Range("E73").Value = 4
$ws.Range("F73").Value = "Fruta"
$ws.Range("G73").Value = 100103
$ws.Range("H73").Value = "Frutos de hueso (carozo)"
$ws.Range("I73").Value = 100103002
$ws.Range("J73").Value = "Ciruela"
$ws.Range("K73").Value = "Angeleno"
$ws.Range("L73").Value = "Segunda"
$ws.Range("M73").Value = 16
$ws.Range("N73").Value = 170000
$ws.Range("O73").Value = 180000
$ws.Range("P73").Value = 175000
$ws.Range("Q73").Value = "$/bins (450 kilos)"
$ws.Range("R73").Value = "Región de O'Higgins"
$ws.Range("S73").Value = 389
$ws.Range("T73").Value = 450
